$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Units" row label becomes "_units" (a new "_weight" row is being added below it)
$ws.Range("A2").Value = "_units"

# Insert a new row at position 3 (pushes Illite..UGas down by one row)
$ws.Rows.Item(3).Insert()

# Populate the new "_weight" row with per-column log weights
$ws.Range("A3").Value = "_weight"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0.8
$ws.Range("F3").Value = 0.7
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 0.8
$ws.Range("I3").Value = 0.3
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 1

# The AutoFilter range grew by one row because of the inserted row
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Constants!_FilterDatabase" -or $n.Name -eq "_FilterDatabase") {
        $n.RefersTo = "=Constants!`$A`$2:`$F`$9"
    }
}

# Restore the active-cell selection
[void]$ws.Range("D10").Select()
